{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// 1) \"--\" + \"Note which weapon you have chose\" + \"--\" -> single run \"--Note which weapon you have chose--\"\nitems[0].getRange().insertText(\"--Note which weapon you have chose--\", \"Replace\");\n\n// 2)-4) Add a first-line indent (0.5in = 36pt) to the three dialogue paragraphs\n//    that follow the banner paragraph.\nitems[1].firstLineIndent = 36;\nitems[2].firstLineIndent = 36;\nitems[3].firstLineIndent = 36;\n\n// 5) Same first-line indent on the \"You shudder...\" paragraph, and relocate the\n//    \"_GoBack\" bookmark here (Word keeps _GoBack pinned to the last edit point).\nitems[4].firstLineIndent = 36;\n\nawait context.sync();\n\n// Remove the old bookmark (currently sitting in the final, near-empty paragraph)\n// before re-inserting it at the start of paragraph 4.\ncontext.document.deleteBookmark(\"_GoBack\");\nitems[4].getRange(\"Start\").insertBookmark(\"_GoBack\");\n\n// 6) \"--\" + \"A year later\" + \"--\" -> single run \"--A year later--\"\nitems[5].getRange().insertText(\"--A year later--\", \"Replace\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Banner paragraph: \"--\" + \"Note which weapon you have chose\" + \"--\"\n#    -> single run \"--Note which weapon you have chose--\".\n#    Two-step Text= forces the engine to actually rewrite (and merge) the\n#    run even though the concatenated text reads the same either way, and\n#    keeps the original run formatting (Bookerly / bold / color).\n$p1 = $d.Paragraphs.Item(1)\n$r1 = $p1.Range\n$r1.MoveEnd(1, -1) | Out-Null\n$r1.Text = \"placeholder\"\n$r1b = $p1.Range\n$r1b.MoveEnd(1, -1) | Out-Null\n$r1b.Text = \"--Note which weapon you have chose--\"\n\n# 2)-4) Add a first-line indent (0.5in = 36pt) to the three dialogue\n#    paragraphs that follow the banner paragraph.\n$d.Paragraphs.Item(2).Range.ParagraphFormat.FirstLineIndent = 36\n$d.Paragraphs.Item(3).Range.ParagraphFormat.FirstLineIndent = 36\n$d.Paragraphs.Item(4).Range.ParagraphFormat.FirstLineIndent = 36\n\n# 5) Same first-line indent on the \"You shudder...\" paragraph, and move the\n#    \"_GoBack\" bookmark here from the final (near-empty) paragraph.\n$d.Paragraphs.Item(5).Range.ParagraphFormat.FirstLineIndent = 36\n\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n$r5 = $d.Paragraphs.Item(5).Range\n$r5.Collapse(1) | Out-Null\n$d.Bookmarks.Add(\"_GoBack\", $r5)\n\n# 6) Banner paragraph: \"--\" + \"A year later\" + \"--\" -> single run \"--A year later--\"\n$p6 = $d.Paragraphs.Item(6)\n$r6 = $p6.Range\n$r6.MoveEnd(1, -1) | Out-Null\n$r6.Text = \"placeholder\"\n$r6b = $p6.Range\n$r6b.MoveEnd(1, -1) | Out-Null\n$r6b.Text = \"--A year later--\"\n"}
